$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "AuthActivity.kt"
$ws.Range("A2").Value = "binding"
